$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row at 9 (this becomes the "Ie" row) - shifts old rows 9-26
# down to 10-27, and Excel auto-adjusts all formulas / shared formulas /
# merged cells referencing those rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(9).Insert()

# ---------------------------------------------------------------------------
# Row 3 (VC): add measured column D
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = 3.124

# ---------------------------------------------------------------------------
# Row 4 (VB): add measured column D, and fill in AC analysis Rin/Rout block
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = 0.704

$ws.Range("H4").Value = 0.006
$ws.Range("I4").Formula = "=3.796*10^-6"
$ws.Range("L4").Value = 0.1
$ws.Range("M4").Formula = "=33.468*10^-6"

# ---------------------------------------------------------------------------
# Row 5 (VE): add measured column D
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = 0.028

# ---------------------------------------------------------------------------
# Row 7 (Ib): add measured column D (scientific format)
# ---------------------------------------------------------------------------
$ws.Range("D7").Value = 0.00028129499999999998
$ws.Range("D7").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# Row 8 (Ic): add measured column D (scientific format)
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = 0.00191
$ws.Range("D8").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# Row 10 (was 9, IR1): add measured column D (scientific format)
# ---------------------------------------------------------------------------
$ws.Range("D10").Value = 0.0002804
$ws.Range("D10").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# Row 11 (was 10, IR2): add measured column D (scientific format)
# ---------------------------------------------------------------------------
$ws.Range("D11").Value = 0.000266007
$ws.Range("D11").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# New string cells, set in the order that matches the shared-string table
# of the target workbook: H21, D15, D16, L13, A35, A36, B36, B35, A37, A38,
# B37, A20, D20, A9
# ---------------------------------------------------------------------------
$ws.Range("H21").Value = "7.6k"
$ws.Range("D15").Value = "9.7k"
$ws.Range("D16").Value = "375K"
$ws.Range("L13").Value = "``"
$ws.Range("A35").Value = "R1"
$ws.Range("A36").Value = "R2"
$ws.Range("B36").Value = "3.108k"
$ws.Range("B35").Value = "47.405k"
$ws.Range("A37").Value = "RC"
$ws.Range("A38").Value = "RE"
$ws.Range("B37").Value = "9.439k"
$ws.Range("A20").Value = "Swing Assymetry Voltage"
$ws.Range("D20").Value = "80m"
$ws.Range("A9").Value = "Ie"

# ---------------------------------------------------------------------------
# Row 9 (new, Ie): sums of Ib/Ic
# ---------------------------------------------------------------------------
$ws.Range("C9").Formula = "=SUM(C7,C8)"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Formula = "=SUM(D7,D8)"
$ws.Range("D9").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# Row 14 (was 13, F-3dB Dominant): add measured column D
# ---------------------------------------------------------------------------
$ws.Range("D14").Value = 800

# ---------------------------------------------------------------------------
# Row 38 (RE): measured value
# ---------------------------------------------------------------------------
$ws.Range("B38").Value = 164.83

# ---------------------------------------------------------------------------
# Column widths / selection cosmetics
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(9).AutoFit()
$ws.Columns.Item(13).AutoFit()

$ws.Range("I16").Select()
